$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Brown, et al.")

# Insert a new column before column I (9th column), shifting the existing
# "Amount for 10mL" / "Amount for 5mL" / "Amount added (g)" columns right by
# one and duplicating the formatting of the old column I for the new cells.
$ws.Columns.Item(9).Insert()

# Populate the new column I with the "Amount for 15mL" calculations.
$ws.Cells.Item(4, 9).Value = "Amount for`n15mL"
$ws.Cells.Item(5, 9).Value = "90mg"
$ws.Cells.Item(6, 9).Value = "15mL"

# Extend the print area to cover the new column (K -> L).
$ws.PageSetup.PrintArea = '$A$1:$L$6'

# Match the selection Excel leaves behind after the edit.
$ws.Range("A1:L6").Select() | Out-Null

# Scale the printed page to 96%.
$ws.PageSetup.Zoom = 96
